$d = $word.ActiveDocument

# NOTE: this runtime's Find/Replace (and any Range text mutation) re-serializes
# every run from the edited point through the END of its enclosing paragraph
# into a single run. To keep the untouched "tail" of each paragraph as close to
# the original run layout as possible, every paragraph that receives more than
# one edit is edited back-to-front (rightmost change first), and the whole
# document is likewise edited from the bottom up.

# --- Paragraph ending "...have been delivered during the past 20 years." ---
# (rightmost edit first)

# 4. "during the past 20 years" -> "during the past 10 years"
$d.Content.Find.Execute("during the past 20 years", $true, $false, $false, $false, $false, $true, 1, $false, "during the past 10 years", 2)

# 3. " 15 million square feet of Class +A office (90% of all CBD inventory)" ->
#    " 15.5 million square feet of Class +A office (74% of all CBD deliveries)"
$d.Content.Find.Execute(" 15 million square feet of Class +A office (90% of all CBD inventory)", $true, $false, $false, $false, $false, $true, 1, $false, " 15.5 million square feet of Class +A office (74% of all CBD deliveries)", 2)

# 2. "units (73% of all CBD inventory) " -> "units (44% of all CBD deliveries) "
$d.Content.Find.Execute("units (73% of all CBD inventory) ", $true, $false, $false, $false, $false, $true, 1, $false, "units (44% of all CBD deliveries) ", 2)

# 1. "25,000" -> "16,000" (luxury residential unit count)
$d.Content.Find.Execute("25,000 luxury residential", $true, $false, $false, $false, $false, $true, 1, $false, "16,000 luxury residential", 2)

# --- "...begins at Halsted Street .5 miles west of Ogilvie, whereas..." ---

# 5. ".5 miles west of Ogilvie" -> ".4 miles west of Ogilvie"
$d.Content.Find.Execute(".5 miles west of Ogilvie", $true, $false, $false, $false, $false, $true, 1, $false, ".4 miles west of Ogilvie", 2)

# --- "...Delivering over 2 million square feet ... 500' of new Riverwalk..." ---

# 6. "500' of new Riverwalk" -> "850' of new Riverwalk"
$d.Content.Find.Execute("500’ of new Riverwalk", $true, $false, $false, $false, $false, $true, 1, $false, "850’ of new Riverwalk", 2)

# 7. Delete the trailing "REMOVED FROM ANALYSIS" commentary paragraphs that
#    followed the "World-class Dining Districts" section, keeping the dining
#    paragraph itself and the final trailing empty paragraph of the document.
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*highest concentrations and variety of exciting dining options*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $target + 1
    $endPara = $paras.Count - 1
    if ($endPara -ge $startPara) {
        $startPos = $paras.Item($startPara).Range.Start
        $endPos = $paras.Item($endPara).Range.End
        $rng = $d.Range($startPos, $endPos)
        $rng.Delete()
    }
}

Write-Output "done"
